$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) mirrors the existing header styling (bold, centered,
# thin-bordered like B1:G1) by copying G1's format into H1, then overwriting
# the copied value with the new header text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data rows: every pitcher appearance recorded this save (value 1).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
